$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.2366412213740458
    "C2" = 0.5076335877862596
    "J2" = 0.01145038167938931
    "P2" = 0.1717557251908397
    "S2" = 0.07251908396946564
    "C3" = 0.02238805970149254
    "J3" = 0.04477611940298507
    "P3" = 0.7014925373134329
    "S3" = 0.2313432835820896
    "J4" = 0.02631578947368421
    "P4" = 0.6052631578947368
    "S4" = 0.3684210526315789
    "B6" = 0.06842105263157895
    "D6" = 0.02631578947368421
    "F6" = 0.07368421052631578
    "J6" = 0.1736842105263158
    "O6" = 0.04210526315789474
    "Q6" = 0.1894736842105263
    "R6" = 0.07368421052631578
    "S6" = 0.3526315789473684
    "B7" = 0.0975609756097561
    "D7" = 0.01951219512195122
    "F7" = 0.04878048780487805
    "J7" = 0.1024390243902439
    "O7" = 0.03414634146341464
    "Q7" = 0.2048780487804878
    "R7" = 0.08780487804878048
    "S7" = 0.4048780487804878
    "B8" = 0.09009009009009009
    "D8" = 0.002252252252252252
    "F8" = 0.07207207207207207
    "J8" = 0.09234234234234234
    "O8" = 0.04054054054054054
    "Q8" = 0.2252252252252252
    "R8" = 0.07657657657657657
    "S8" = 0.4009009009009009
    "B9" = 0.08609271523178808
    "D9" = 0.01986754966887417
    "F9" = 0.09933774834437085
    "J9" = 0.09271523178807947
    "O9" = 0.03973509933774835
    "Q9" = 0.2384105960264901
    "R9" = 0.06622516556291391
    "S9" = 0.3576158940397351
    "B10" = 0.1040299906279288
    "D10" = 0.02624179943767573
    "F10" = 0.05716963448922212
    "J10" = 0.1349578256794752
    "O10" = 0.03561387066541706
    "Q10" = 0.2567947516401125
    "R10" = 0.05998125585754452
    "S10" = 0.3252108716026242
    "G11" = 0.1226765799256506
    "J11" = 0.08550185873605948
    "K11" = 0.1561338289962825
    "L11" = 0.6282527881040892
    "S11" = 0.007434944237918215
    "G12" = 0.8255813953488372
    "J12" = 0.1337209302325581
    "L12" = 0.005813953488372093
    "S12" = 0.03488372093023256
    "G13" = 0.8095238095238095
    "J13" = 0.1904761904761905
    "F15" = 0.01234567901234568
    "H15" = 0.1440329218106996
    "I15" = 0.0411522633744856
    "J15" = 0.3127572016460906
    "K15" = 0.08230452674897119
    "M15" = 0.0205761316872428
    "O15" = 0.09876543209876543
    "S15" = 0.2880658436213992
    "F16" = 0.01333333333333333
    "H16" = 0.2333333333333333
    "I16" = 0.05333333333333334
    "J16" = 0.4133333333333333
    "K16" = 0.1266666666666667
    "M16" = 0.02666666666666667
    "O16" = 0.04666666666666667
    "S16" = 0.08666666666666667
    "F17" = 0.02070393374741201
    "H17" = 0.1884057971014493
    "I17" = 0.08074534161490683
    "J17" = 0.4347826086956522
    "K17" = 0.1118012422360248
    "M17" = 0.01863354037267081
    "O17" = 0.06625258799171843
    "S17" = 0.07867494824016563
    "F18" = 0.02173913043478261
    "H18" = 0.2391304347826087
    "I18" = 0.08695652173913043
    "J18" = 0.3840579710144927
    "K18" = 0.08695652173913043
    "M18" = 0.02898550724637681
    "O18" = 0.05797101449275362
    "S18" = 0.09420289855072464
    "F19" = 0.01771653543307087
    "H19" = 0.2421259842519685
    "I19" = 0.08070866141732283
    "J19" = 0.3671259842519685
    "K19" = 0.1151574803149606
    "M19" = 0.02066929133858268
    "N19" = 0.001968503937007874
    "O19" = 0.06791338582677166
    "S19" = 0.08661417322834646
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

Write-Host "Updated $($values.Count) cells"